$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update Estados Unidos (row 4) totals ---
$ws.Range("B4").Value = 1359348
$ws.Range("C4").Value = 12039
$ws.Range("E4").Value = 1038217
$ws.Range("G4").Value = 478
$ws.Range("H4").Value = 80515

# --- Barein overtakes Moldavia in the ranking (rows 60/61 swap) ---
# Row 60 becomes Barein with refreshed stats; row 61 becomes Moldavia keeping its old stats.
$ws.Range("A60").Value = "Barein"
$ws.Range("B60").Value = 4941
$ws.Range("C60").Value = 167
$ws.Range("D60").Value = 2070
$ws.Range("E60").Value = 2863
$ws.Range("F60").Value = 2
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 8

$ws.Range("A61").Value = "Moldavia"
$ws.Range("B61").Value = 4927
$ws.Range("C61").Value = 60
$ws.Range("D61").Value = 1958
$ws.Range("E61").Value = 2800
$ws.Range("F61").Value = 237
$ws.Range("G61").Value = 8
$ws.Range("H61").Value = 169

# --- Update Senegal (row 84) ---
$ws.Range("E84").Value = 1040
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = 19

# --- Update Costa Rica (row 108) ---
$ws.Range("B108").Value = 792
$ws.Range("C108").Value = 12
$ws.Range("D108").Value = 501
$ws.Range("E108").Value = 285
